# cs5010_Dict_DRAFTv2_100720.xlsx - "inflation 1995 ipynb work"
# Populate the three previously-blank data-dictionary rows (30-32) describing
# the new "amount" / "inflation rate" / "cumulative" columns added for the
# 1995-dollars inflation work, add their shared hyperlink, and re-select the
# last-edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order A..AD matches the existing header row (row 4):
# db_Col, db_ColName, Description, Data_Type, Quant/Categorical/Filter, #Values,
# Units, IsAnnual, IsQuarterly, IsMonthly, IsWeekly, Min_Date, Max_Date, IsUS,
# IsCensusDivision, IsState, IsCounty, IsCity, IsZip, IsGEOID, IsMSA,
# Spatial_Identifier, Min_Value, Max_Value, Source, Source_Type,
# Source_Location, Source_Name, Source_Col_Name, Calculated

$rows = @{
    30 = @(24, "amount", "Dollar amount in 1995 dollars", "float64", "Quant", "Continuous", "Dollars",
           1, 0, 0, 0, 1995, 2020, 1, 0, 0, 0, 0, 0, 0, 0, "US", 1, 1.71,
           "https://www.in2013dollars.com/us/inflation/1995", "web scraped to CSV", "../DataSet/",
           "in2013Dollars.com from Bureau Of Labor Statistics", "amount", "No")
    31 = @(25, "inflation rate", "year over year inflation rate", "float64", "Quant", "Continuous", "Percentage",
           1, 0, 0, 0, 1995, 2020, 1, 0, 0, 0, 0, 0, 0, 0, "US", 0, 0.0384,
           "https://www.in2013dollars.com/us/inflation/1995", "web scraped to CSV", "../DataSet/",
           "in2013Dollars.com from Bureau Of Labor Statistics", "inflation rate", "No")
    32 = @(26, "cumulative", "cumulative from 1995 inflation rate", "float64", "Quant", "Continuous", "Percentage",
           1, 0, 0, 0, 1995, 2020, 1, 0, 0, 0, 0, 0, 0, 0, "US", 0, 0.7079,
           "https://www.in2013dollars.com/us/inflation/1995", "web scraped to CSV", "../DataSet/",
           "in2013Dollars.com from Bureau Of Labor Statistics", "cumulative", "No")
}

foreach ($r in 30, 31, 32) {
    $vals = $rows[$r]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# The Source column (Y) for each new row is a hyperlink to the in2013dollars
# inflation page; display text is simply the URL itself. These three cells
# get the borderless "Hyperlink" look (unlike the earlier Y16/Y17 links,
# which kept the table's cell border), so strip the inherited border first.
$url = "https://www.in2013dollars.com/us/inflation/1995"
foreach ($addr in "Y30", "Y31", "Y32") {
    $cell = $ws.Range($addr)
    for ($edge = 1; $edge -le 4; $edge++) {
        $cell.Borders.Item($edge).LineStyle = -4142   # xlLineStyleNone
    }
    $ws.Hyperlinks.Add($cell, $url) | Out-Null
}

# Leave the selection where the author finished editing.
$ws.Range("AD31").Select() | Out-Null
